# Add 32 new "scrambled option" translations as two extra columns (L and P)
# of Korean sentence variants, one pair of new shared strings per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "통에 모든 해적이 기대어 있다"
$ws.Range("P2").Value = "통 하나에 모든 해적이 기대어 있다"
$ws.Range("L3").Value = "모든 통에 해적 한사람이 기대어 있다"
$ws.Range("P3").Value = "모든 통에 해적이 기대어 있다"
$ws.Range("L4").Value = "물고기 한마리를 해적 한사람이 잡았다"
$ws.Range("P4").Value = "물고기를 모든 해적이 잡았다"
$ws.Range("L5").Value = "모든 물고기를 해적 한사람이 잡았다"
$ws.Range("P5").Value = "모든 물고기를 해적이 잡았다"
$ws.Range("L6").Value = "낚시대 하나를 모든 해적이 잡고 있다"
$ws.Range("P6").Value = "낚시대를 모든 해적이 잡고 있다"
$ws.Range("L7").Value = "모든 낚시대를 해적 한사람이 잡고 있다"
$ws.Range("P7").Value = "모든 낚시대를 해적이 잡고 있다"
$ws.Range("L8").Value = "상어 한마리를 모든 해적이 먹였다"
$ws.Range("P8").Value = "상어를 모든 해적이 먹였다"
$ws.Range("L9").Value = "모든 상어를 해적 한사람이 먹였다"
$ws.Range("P9").Value = "모든 상어를 해적이 먹였다"
$ws.Range("L10").Value = "병 하나를 모든 해적이 잡고 있다"
$ws.Range("P10").Value = "병을 모든 해적이 잡고 있다"
$ws.Range("L11").Value = "모든 병을 해적 한사람이 잡고 있다"
$ws.Range("P11").Value = "모든 병을 해적이 잡고 있다"
$ws.Range("L12").Value = "물고기 한마리를 모든 상어가 물고 있다"
$ws.Range("P12").Value = "물고기를 모든 상어가 물고 있다"
$ws.Range("L13").Value = "모든 물고기를 상어 한마리가 물고 있다"
$ws.Range("P13").Value = "모든 물고기를 상어가 물고 있다"
$ws.Range("L14").Value = "해적 한사람을 모든 상어가 공격했다"
$ws.Range("P14").Value = "해적을 모든 상어가 공격했다"
$ws.Range("L15").Value = "모든 해적을 상어 한마리가 공격했다"
$ws.Range("P15").Value = "모든 해적을 상어가 공격했다"
$ws.Range("L16").Value = "강아지 한마리를 모든 여자아이가 쓰다듬고 있다"
$ws.Range("P16").Value = "강아지를 모든 여자아이가 쓰다듬고 있다"
$ws.Range("L17").Value = "모든 강아지를 여자아이 하나가 쓰다듬고 있었다"
$ws.Range("P17").Value = "모든 강아지를 여자아이가 쓰다듬고 있었다"

# Widen column G to fit the longer Korean text that was already there
$ws.Columns(7).ColumnWidth = 26.3

# Match the saved view: zoomed in a bit, scrolled so column B is at the
# left edge, with the selection resting one row below the new data.
$win = $excel.ActiveWindow
$win.Zoom = 116
$ws.Range("P18").Select()
